$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '63.241.12'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.98%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.569.23'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +1.30%  '
$ws.Range("E4").Value = '  -0.02%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '584.73'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +3.43%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '148.22'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +2.14%  '
$ws.Range("E7").Value = '  -0.03%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.602'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +3.70%  '
$ws.Range("E9").Value = '  +4.22%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '5.64'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +0.83%  '
$ws.Range("E11").Value = '  +0.46%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.357'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +1.86%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '27.51'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +2.22%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '3.027.71'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +1.13%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '63.179.23'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +0.84%  '
$ws.Range("E16").Value = '  +4.87%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '2.570.63'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("E18").Value = '  -0.24%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '342.63'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +2.89%  '
$ws.Range("E20").Value = '  +3.81%  '
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("E22").Value = '  +0.05%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '66.69'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +3.47%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '2.693.06'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +1.59%  '
$ws.Range("E25").Value = '  +3.51%  '
$ws.Range("E26").Value = '  +1.49%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '8.22'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +14.00%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '8.49'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +2.54%  '
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("E30").Value = '  -0.31%  '
$ws.Range("E31").Value = '  +8.09%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0825'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +2.87%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '459.34'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +13.65%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.64'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +4.50%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '176.86'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +0.27%  '
$ws.Range("E36").Value = '  +2.76%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '19.25'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +2.32%  '
$ws.Range("E38").Value = '  +4.76%  '
$ws.Range("E40").Value = '  +0.36%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '151.23'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -0.85%  '
$ws.Range("E43").Value = '  +2.68%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '21.10'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +2.79%  '
$ws.Range("E45").Value = '  +7.27%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.614'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +2.50%  '
$ws.Range("E47").Value = '  +2.83%  '
$ws.Range("E48").Value = '  +2.62%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '18.41'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +1.70%  '
$ws.Range("E50").Value = '  +0.21%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '11.39'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.03%  '
